# Change the year in the astromap link: 2018 -> 2022, and collapse the
# paragraph's four differently-formatted runs into the plain-text run
# structure used by the edited document (an empty run followed by a
# single unformatted run holding the whole sentence).

$d = $word.ActiveDocument

# Locate the paragraph that holds the astromap credit/link sentence.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Schaubilder*CzechGlobe*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range
    # Exclude the trailing paragraph mark from the range we replace.
    $r.MoveEnd(1, -1) | Out-Null
    $r.Delete() | Out-Null
    $r.Collapse(1) | Out-Null

    $newSentence = "Die Schaubilder in diesem Dokument wurden von Jan Hollan, CzechGlobe, bereitgestellt. (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
           '<w:body>' +
           '<w:p w14:paraId="57F3A1B8" w14:textId="5ACF7167" w:rsidR="00550C9F" w:rsidRPr="0065155F" w:rsidRDefault="00550C9F" w:rsidP="0065155F">' +
           '<w:pPr>' +
           '<w:pStyle w:val="BasicParagraph"/>' +
           '<w:pBdr>' +
           '<w:top w:val="single" w:sz="4" w:space="5" w:color="auto"/>' +
           '<w:left w:val="single" w:sz="4" w:space="4" w:color="auto"/>' +
           '<w:bottom w:val="single" w:sz="4" w:space="10" w:color="auto"/>' +
           '<w:right w:val="single" w:sz="4" w:space="4" w:color="auto"/>' +
           '</w:pBdr>' +
           '<w:spacing w:line="240" w:lineRule="auto"/>' +
           '<w:ind w:right="-90"/>' +
           '<w:jc w:val="center"/>' +
           '<w:rPr>' +
           '<w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>' +
           '<w:sz w:val="19"/>' +
           '<w:lang w:val="de-DE"/>' +
           '</w:rPr>' +
           '</w:pPr>' +
           '<w:r/>' +
           '<w:r><w:t>' + $newSentence + '</w:t></w:r>' +
           '</w:p>' +
           '</w:body>' +
           '</w:document>' +
           '</pkg:xmlData>' +
           '</pkg:part>' +
           '</pkg:package>'

    $r.InsertXML($xml) | Out-Null
}
